$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for price cells (column D) so numeric-looking strings
# like "61.936.30" or "1.00" are preserved exactly as text.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D12", "D13", "D15", "D16", "D17", "D19", "D20", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }

# Apply updated values
$ws.Range("D2").Value = "61.936.30"
$ws.Range("D3").Value = "2.416.88"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "562.64"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").Value = "142.97"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").Value = "0.349"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "25.69"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "2.852.18"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "61.894.29"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "2.413.82"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").Value = "322.81"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "6.83"
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "66.53"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").Value = "8.72"
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("D26").Value = "558.77"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "2.536.21"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "0.0₃0934"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "8.17"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("D32").Value = "0.147"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "1.86"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -3.44%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "4.75"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").Value = "153.82"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("D40").Value = "18.49"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "0.995"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").Value = "2.25"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "147.11"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "3.62"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "0.0525"
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.592"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "19.77"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").Value = "0.0920"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "0.0226"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  +0.81%  "
